$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C column) date for rows 2 through 11
# from 2023-09-01 (45170) to 2023-09-05 (45174)
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
